# Generate Report for Handoff
# Adds two newly-handed-off files (43b87450-1188-4359-bc8d-7e05de8a463e.md and
# 92527a98-232b-4d65-aa57-bfe41d66a16d.md) as new rows on all three sheets:
#   Overview (row4/row5), zh-cn (row4/row5), de-de (row4/row5)

$wb = $excel.ActiveWorkbook

$file1 = "43b87450-1188-4359-bc8d-7e05de8a463e.md"
$file2 = "92527a98-232b-4d65-aa57-bfe41d66a16d.md"

$hash1 = "4714b7f7ea52696c766485f5d705c10dafc13b3a"
$hash2 = "d8264ec67ceae6e24a5358a146a60e7057a41c36"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $file1
$wsOverview.Range("B4").Value = "e2e\$file1"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-04 14:44:08"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("A5").Value = $file2
$wsOverview.Range("B5").Value = "e2e\$file2"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-04 14:44:08"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/$file1", $null, $null, "e2e\$file1") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/$file2", $null, $null, "e2e\$file2") | Out-Null

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $file1
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = "$file1.$hash1.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-09-04 14:43:59"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

$wsZh.Range("A5").Value = $file2
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = "$file2.$hash2.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-09-04 14:43:59"
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "False"
$wsZh.Range("P5").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/$file1", $null, $null, $file1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/$file2", $null, $null, $file2) | Out-Null

$wsZh.ListObjects.Item(1).Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $file1
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = "$file1.$hash1.de-de.xlf"
$wsDe.Range("H4").Value = "2016-09-04 14:44:08"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

$wsDe.Range("A5").Value = $file2
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = "$file2.$hash2.de-de.xlf"
$wsDe.Range("H5").Value = "2016-09-04 14:44:08"
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "False"
$wsDe.Range("P5").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/$file1", $null, $null, $file1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/$file2", $null, $null, $file2) | Out-Null

$wsDe.ListObjects.Item(1).Resize($wsDe.Range("A1:P5"))
